$d = $word.ActiveDocument

$replacements = @(
    @("31÷6=", "45÷4="),
    @("57÷2=", "67÷3="),
    @("43÷5=", "40÷5="),
    @("93÷7=", "46÷4="),
    @("14÷6=", "23÷3="),
    @("97÷3=", "97÷9="),
    @("60÷7=", "35÷5="),
    @("44÷2=", "23÷3="),
    @("39÷6=", "66÷8="),
    @("42÷6=", "43÷4="),
    @("32÷5=", "93÷8="),
    @("95÷4=", "51÷2="),
    @("33÷8=", "50÷3="),
    @("13÷2=", "39÷9="),
    @("12÷6=", "67÷9="),
    @("47÷2=", "32÷7="),
    @("71÷8=", "94÷2="),
    @("65÷9=", "23÷4="),
    @("38÷5=", "72÷6="),
    @("21÷7=", "69÷9="),
    @("65÷2=", "66÷6="),
    @("20÷9=", "56÷6="),
    @("23÷6=", "44÷4="),
    @("59÷6=", "70÷7="),
    @("90÷7=", "75÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
